$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.990638256072998
$ws.Range("B1").Value = 2.312676906585693
$ws.Range("C1").Value = 5.051010608673096
$ws.Range("D1").Value = 2.474934339523315
$ws.Range("E1").Value = 1.369349956512451
